$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 11 new rows before row 153, shifting the old
# "Development type" (153-161) and "Voluntary agreement" (162-163)
# modules down to rows 164-174.
$ws.Rows("153:163").Insert()

# Step 2: widen column G from 27 to 29 characters.
$ws.Columns("G").ColumnWidth = 28.17

# Step 3: populate the new rows 153-163 with the "Oil and gas permission
# types" module content.

# Row 153
$ws.Cells.Item(153, 1).Value = "Oil and gas permission types"
$ws.Cells.Item(153, 2).Value = "oilgas-permission-type"
$ws.Cells.Item(153, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(153, 4).Value = "oilgas-permission-types"
$ws.Cells.Item(153, 5).Value = "Oil and gas permission types[]"
$ws.Cells.Item(153, 12).Value = "List of permission types being applied for"
$ws.Cells.Item(153, 13).Value = "enum"
$ws.Cells.Item(153, 14).Value = "MUST"

# Row 154
$ws.Cells.Item(154, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(154, 4).Value = "related-permissions"
$ws.Cells.Item(154, 5).Value = "Related permissions[]"
$ws.Cells.Item(154, 6).Value = "reference"
$ws.Cells.Item(154, 7).Value = "Reference"
$ws.Cells.Item(154, 12).Value = "The reference for the related application that permission was received for"
$ws.Cells.Item(154, 13).Value = "string"
$ws.Cells.Item(154, 14).Value = "MUST"

# Row 155
$ws.Cells.Item(155, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(155, 4).Value = "related-permissions"
$ws.Cells.Item(155, 5).Value = "Related permissions[]"
$ws.Cells.Item(155, 6).Value = "oilgas-permission-type"
$ws.Cells.Item(155, 7).Value = "Oil and gas permission type"
$ws.Cells.Item(155, 12).Value = "An oil and gas related permission type"
$ws.Cells.Item(155, 13).Value = "enum"
$ws.Cells.Item(155, 14).Value = "MUST"

# Row 156
$ws.Cells.Item(156, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(156, 4).Value = "related-permissions"
$ws.Cells.Item(156, 5).Value = "Related permissions[]"
$ws.Cells.Item(156, 6).Value = "decision-date"
$ws.Cells.Item(156, 7).Value = "Decision date"
$ws.Cells.Item(156, 12).Value = "The date when the decision was made, in YYYY-MM-DD format"
$ws.Cells.Item(156, 13).Value = "string"
$ws.Cells.Item(156, 14).Value = "MUST"

# Row 157
$ws.Cells.Item(157, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(157, 4).Value = "related-permissions"
$ws.Cells.Item(157, 5).Value = "Related permissions[]"
$ws.Cells.Item(157, 6).Value = "condition-number"
$ws.Cells.Item(157, 7).Value = "Condition number"
$ws.Cells.Item(157, 12).Value = "Number of any condition being breached"
$ws.Cells.Item(157, 13).Value = "string"
$ws.Cells.Item(157, 14).Value = "MAY"

# Row 158
$ws.Cells.Item(158, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(158, 4).Value = "other-details"
$ws.Cells.Item(158, 5).Value = "Other details"
$ws.Cells.Item(158, 12).Value = "Explanation if other ground is selected"
$ws.Cells.Item(158, 13).Value = "string"
$ws.Cells.Item(158, 14).Value = "MAY"

# Row 159
$ws.Cells.Item(159, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(159, 4).Value = "will-consolidate-permissions"
$ws.Cells.Item(159, 5).Value = "Will consolidate permissions"
$ws.Cells.Item(159, 12).Value = "Is the applicant looking to consolidate permissions?"
$ws.Cells.Item(159, 13).Value = "boolean"
$ws.Cells.Item(159, 14).Value = "MUST"

# Row 160
$ws.Cells.Item(160, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(160, 4).Value = "details"
$ws.Cells.Item(160, 5).Value = "Details"
$ws.Cells.Item(160, 12).Value = "Details about the consolidation or update of permissions"
$ws.Cells.Item(160, 13).Value = "string"
$ws.Cells.Item(160, 14).Value = "MAY"

# Row 161
$ws.Cells.Item(161, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(161, 4).Value = "related-proposals"
$ws.Cells.Item(161, 5).Value = "Related proposals[]"
$ws.Cells.Item(161, 6).Value = "reference"
$ws.Cells.Item(161, 7).Value = "Reference"
$ws.Cells.Item(161, 12).Value = "The reference for the related application"
$ws.Cells.Item(161, 13).Value = "string"
$ws.Cells.Item(161, 14).Value = "MUST"

# Row 162
$ws.Cells.Item(162, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(162, 4).Value = "related-proposals"
$ws.Cells.Item(162, 5).Value = "Related proposals[]"
$ws.Cells.Item(162, 6).Value = "application-type"
$ws.Cells.Item(162, 7).Value = "Application type"
$ws.Cells.Item(162, 12).Value = "The type of planning application"
$ws.Cells.Item(162, 13).Value = "enum"
$ws.Cells.Item(162, 14).Value = "MUST"

# Row 163
$ws.Cells.Item(163, 3).Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Cells.Item(163, 4).Value = "related-proposals"
$ws.Cells.Item(163, 5).Value = "Related proposals[]"
$ws.Cells.Item(163, 6).Value = "decision-date"
$ws.Cells.Item(163, 7).Value = "Decision date"
$ws.Cells.Item(163, 12).Value = "The date when the decision was made, in YYYY-MM-DD format"
$ws.Cells.Item(163, 13).Value = "string"
$ws.Cells.Item(163, 14).Value = "MUST"

# Re-fit the row heights for the new rows (the C column text ends in a
# newline, which would otherwise leave a stray explicit row height behind).
$ws.Rows("153:163").AutoFit()

# Step 4: merge the A and B label columns across the new module's rows.
$ws.Range("A153:A163").Merge()
$ws.Range("B153:B163").Merge()

Write-Host "Edit complete"
